$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" (column G) for both rows
$wsOverview.Range("G2").Value = "2016-09-05 12:04:59"
$wsOverview.Range("G3").Value = "2016-09-05 12:04:59"

# zh-cn sheet: Priority (E) ht -> mt
$wsZhCn.Range("E2").Value = "mt"
$wsZhCn.Range("E3").Value = "mt"

# zh-cn sheet: Correspond Handoff Datetime (H)
$wsZhCn.Range("H2").Value = "2016-09-05 12:04:48"
$wsZhCn.Range("H3").Value = "2016-09-05 12:04:48"

# zh-cn sheet: Correspond Handback DateTime (K)
$wsZhCn.Range("K2").Value = "2016-09-05 12:05:37"
$wsZhCn.Range("K3").Value = "2016-09-05 12:05:37"

# de-de sheet: Priority (E) ht -> mt (shares the same shared-string as zh-cn!E2/E3)
$wsDeDe.Range("E2").Value = "mt"
$wsDeDe.Range("E3").Value = "mt"

# de-de sheet: Correspond Handback DateTime (K)
$wsDeDe.Range("K2").Value = "2016-09-05 12:05:56"
$wsDeDe.Range("K3").Value = "2016-09-05 12:05:56"

# de-de sheet: Correspond Handoff Datetime (H) -> uses same shared string as Overview!G2/G3
$wsDeDe.Range("H2").Value = "2016-09-05 12:04:59"
$wsDeDe.Range("H3").Value = "2016-09-05 12:04:59"
